# Add a new "2021" column (O) to the wastewater-treated indicator table,
# mirroring the formatting already used by the adjacent "2020" column (N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column N (rows 3-14, the only rows that have
# data in column N) into column O before writing the new values, so O
# inherits the same borders/number formats/fonts as N.
$ws.Range("N3:N14").Copy()
$ws.Range("O3:O14").PasteSpecial(-4122)  # xlPasteFormats

# New 2021 figures for each oblast / the republic as a whole.
$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 97
$ws.Range("O6").Value = 96.2
$ws.Range("O7").Value = 62.7
$ws.Range("O8").Value = 100
$ws.Range("O9").Value = 100
$ws.Range("O10").Value = "-"
$ws.Range("O11").Value = 100
$ws.Range("O12").Value = 57.9
$ws.Range("O13").Value = 100
$ws.Range("O14").Value = "-"

# Match the author's final selection state.
$ws.Range("O17").Select() | Out-Null
